# Apply weekly fruit/vegetable price updates (Higo - Vega Central Mapocho de Santiago)
# This swaps the Fecha/Volumen/Precio/Origen data between matching row pairs
# (2<->12, 3<->13, 6<->14, 7<->15) per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$rowA,
        [int]$rowB
    )

    $colsToSwap = @("D", "M", "N", "O", "P", "R", "S")

    foreach ($col in $colsToSwap) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

Swap-RowData 2  12
Swap-RowData 3  13
Swap-RowData 6  14
Swap-RowData 7  15
